# fix(read from pdf): split columns and split cells in rows
#
# The "Name" column actually held "Forename Surname" pairs, and the
# "Date Of Birth" header was too verbose. This:
#   1. Splits the Name column into separate Name / Surname columns
#      (inserts a new column B and moves the surname into it).
#   2. Splits each "Forename Surname" data cell into its two halves.
#   3. Shortens the "Date Of Birth" header to "D.O.B".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before B - this pushes Age/Sex/Colour/Height/
#    Date Of Birth/Wanted? one column to the right (B:G -> C:H) and
#    carries each column's formatting (incl. the date & boolean number
#    formats on what are now G and H) along with it automatically.
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").ColumnWidth = 8

# 2. Header row: label the new column and tighten the D.O.B header
#    (now in column G after the insert above).
$ws.Range("B1").Value = "Surname"
$ws.Range("G1").Value = "D.O.B"

# 3. Split the "Forename Surname" values in column A into Name (A) /
#    Surname (B) for each data row.
$ws.Range("A2").Value = "Adam"
$ws.Range("B2").Value = "Smith"

$ws.Range("A3").Value = "Janet"
$ws.Range("B3").Value = "Jones"

# 4. Match the cursor position left in the source workbook.
[void]$ws.Range("G1").Select()
